$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.141.15'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +5.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.20'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +2.59%  '

$ws.Range("E4").Value = '  -0.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.22'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +4.73%  '

$ws.Range("E6").Value = '  -0.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5221'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +3.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4090'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +4.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08517'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +1.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.86'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.40'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +9.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.428'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +3.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.922.35'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +3.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.409'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +1.87%  '

$ws.Range("E16").Value = '  -0.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.76'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +5.02%  '

$ws.Range("E18").Value = '  +1.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06683'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.43'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +3.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -0.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.019'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.149.33'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +5.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.34'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +2.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.206'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.139.01'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +2.69%  '

$ws.Range("E27").Value = '  +2.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.85'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +0.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.448'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +0.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.45'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +1.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.085'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +3.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1060'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +2.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.063'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +5.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.641'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02493'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +1.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06629'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2215'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +2.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.238'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +4.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.195'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +3.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.901'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6563'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +2.87%  '

$ws.Range("E42").Value = '  +1.23%  '

$ws.Range("E43").Value = '  +4.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6168'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +2.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.33'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +2.39%  '

$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.083'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +3.74%  '

$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.74'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +1.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.173'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +10.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.81'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +4.17%  '
